$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 206
$ws.Range("B9").Value = "Reverse Linked List"
$ws.Range("D9").Value = "Walk and push to a stack, pop the stack"
$ws.Range("E9").Value = "move 3 ptrs switch direction or recursive (prev,current)"
$ws.Range("C9").Value = "Stack/3 Ptrs/Recursive"

$ws.Range("A10").Value = 217
$ws.Range("B10").Value = "Contains Duplicate"
$ws.Range("C10").Value = "HashSet/Array.Sort/Old school O(n^2)"

$ws.Range("A11").Value = 219
$ws.Range("B11").Value = "Contains Duplicate 2"
$ws.Range("C11").Value = "HashMap/Dictionary/unordered_map"
$ws.Range("D11").Value = "have a map, loop add key, if key exists check abs value with n, return true, else assign new index to current map value. "

$ws.Range("D14").Select()
